# Auto-generated Excel COM-interop script to apply numeric updates
# described by the OOXML diff for Sheets/Ifrit_Profits.xlsx (FFXIV Leve profit tracker).
# Each worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW) gets a batch of cell value updates;
# a handful of cells are cleared (HQ profit no longer applicable) or newly populated
# (HQ profit now applicable), matching the source diff exactly.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1069.3334
$ws.Range("J17").Value = 1069.3334
$ws.Range("L17").Value = 3208.0002
$ws.Range("N17").Value = -3544.0002
$ws.Range("H19").Value = 681.0909
$ws.Range("I19").Value = 540.2
$ws.Range("J19").Value = 798.5
$ws.Range("K19").Value = 540.2
$ws.Range("L19").Value = 798.5
$ws.Range("M19").Value = -365.2
$ws.Range("N19").Value = -1148.5
$ws.Range("H112").Value = 1377.7368
$ws.Range("J112").Value = 1558.4667
$ws.Range("L112").Value = 4675.4001
$ws.Range("N112").Value = -6891.4001
$ws.Range("H132").Value = 9500
$ws.Range("I132").Value = 11751.875
$ws.Range("J132").Value = 492.5
$ws.Range("K132").Value = 35255.625
$ws.Range("L132").Value = 1477.5
$ws.Range("M132").Value = -32725.625
$ws.Range("N132").Value = -6537.5
$ws.Range("H137").Value = 24393060
$ws.Range("I137").Value = 1712.6875
$ws.Range("K137").Value = 5138.0625
$ws.Range("M137").Value = -2588.0625
$ws.Range("H138").Value = 3298.0356
$ws.Range("I138").Value = 2826.9714
$ws.Range("J138").Value = 4083.1428
$ws.Range("K138").Value = 8480.914199999999
$ws.Range("L138").Value = 12249.4284
$ws.Range("M138").Value = -3340.914199999999
$ws.Range("N138").Value = -22529.4284

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2877654.5
$ws.Range("I132").Value = 6901130
$ws.Range("J132").Value = 3743.7144
$ws.Range("K132").Value = 20703390
$ws.Range("L132").Value = 11231.1432
$ws.Range("M132").Value = -20700860
$ws.Range("N132").Value = -16291.1432

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 50780
$ws.Range("J132").Value = 50780
$ws.Range("L132").Value = 50780
$ws.Range("N132").Value = -60900
$ws.Range("H141").Value = 57813.547
$ws.Range("J141").Value = 53471.11
$ws.Range("L141").Value = 53471.11
$ws.Range("N141").Value = -63831.11

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2076.5652
$ws.Range("I86").Value = 2177.6428
$ws.Range("J86").Value = 1919.3334
$ws.Range("K86").Value = 2177.6428
$ws.Range("L86").Value = 1919.3334
$ws.Range("M86").Value = -1054.6428
$ws.Range("N86").Value = -4165.3334
$ws.Range("H89").Value = 2076.5652
$ws.Range("I89").Value = 2177.6428
$ws.Range("J89").Value = 1919.3334
$ws.Range("K89").Value = 10888.214
$ws.Range("L89").Value = 9596.666999999999
$ws.Range("M89").Value = -5272.214
$ws.Range("N89").Value = -20828.667
$ws.Range("H93").Value = 9751.166999999999
$ws.Range("I93").Value = 7910.364
$ws.Range("K93").Value = 7910.364
$ws.Range("M93").Value = -6038.364
$ws.Range("H94").Value = 1938.4445
$ws.Range("I94").Value = 1101
$ws.Range("J94").Value = 2177.7144
$ws.Range("K94").Value = 1101
$ws.Range("L94").Value = 2177.7144
$ws.Range("M94").Value = -650
$ws.Range("N94").Value = -3079.7144
$ws.Range("H132").Value = 2431.037
$ws.Range("I132").Value = 2194.8696
$ws.Range("J132").Value = 3789
$ws.Range("K132").Value = 6584.6088
$ws.Range("L132").Value = 11367
$ws.Range("M132").Value = -4054.6088
$ws.Range("N132").Value = -16427

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H86").Value = 771
$ws.Range("J86").Value = 919.8
$ws.Range("L86").Value = 2759.4
$ws.Range("N86").Value = -5131.4
$ws.Range("H89").Value = 771
$ws.Range("J89").Value = 919.8
$ws.Range("L89").Value = 8278.199999999999
$ws.Range("N89").Value = -20134.2
$ws.Range("H103").Value = 5667825.5
$ws.Range("J103").Value = 1482
$ws.Range("L103").Value = 4446
$ws.Range("N103").Value = -6204
$ws.Range("H105").Value = 892000000
$ws.Range("J105").Value = 892000000
$ws.Range("L105").Value = 2676000000
$ws.Range("N105").Value = -2676005242
$ws.Range("H107").Value = 154085.16
$ws.Range("I107").Value = 111310.336
$ws.Range("K107").Value = 333931.008
$ws.Range("M107").Value = -332011.008
$ws.Range("H110").Value = 3780
$ws.Range("J110").Value = 3780
$ws.Range("L110").Value = 11340
$ws.Range("N110").Value = -19520
$ws.Range("H113").Value = 636.1111
$ws.Range("I113").Value = 643.8182
$ws.Range("J113").Value = 624
$ws.Range("K113").Value = 1931.4546
$ws.Range("L113").Value = 1872
$ws.Range("M113").Value = 238.5454
$ws.Range("N113").Value = -6212
$ws.Range("H131").Value = 2809.8728
$ws.Range("J131").Value = 1783.2642
$ws.Range("L131").Value = 5349.792600000001
$ws.Range("N131").Value = -15429.7926
$ws.Range("H132").Value = 58824750
$ws.Range("I132").Value = 100001020
$ws.Range("J132").Value = 1511.1428
$ws.Range("K132").Value = 900009180
$ws.Range("L132").Value = 13600.2852
$ws.Range("M132").Value = -900006650
$ws.Range("N132").Value = -18660.2852

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 48496.918
$ws.Range("I80").Value = 2394.5454
$ws.Range("J80").Value = 87506.62
$ws.Range("K80").Value = 2394.5454
$ws.Range("L80").Value = 87506.62
$ws.Range("M80").Value = -1396.5454
$ws.Range("N80").Value = -89502.62
$ws.Range("H83").Value = 48496.918
$ws.Range("I83").Value = 2394.5454
$ws.Range("J83").Value = 87506.62
$ws.Range("K83").Value = 11972.727
$ws.Range("L83").Value = 437533.1
$ws.Range("M83").Value = -6980.726999999999
$ws.Range("N83").Value = -447517.1
$ws.Range("H113").Value = 2542.111
$ws.Range("I113").Value = 3129.8333
$ws.Range("K113").Value = 3129.8333
$ws.Range("M113").Value = -959.8332999999998
$ws.Range("H132").Value = 2537.5833
$ws.Range("I132").Value = 2216.7
$ws.Range("J132").Value = 4142
$ws.Range("K132").Value = 6650.099999999999
$ws.Range("L132").Value = 12426
$ws.Range("M132").Value = -4120.099999999999
$ws.Range("N132").Value = -17486

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2332.6667
$ws.Range("I40").Value = 2290.2222
$ws.Range("J40").Value = 2364.5
$ws.Range("K40").Value = 2290.2222
$ws.Range("L40").Value = 2364.5
$ws.Range("M40").Value = -2154.2222
$ws.Range("N40").Value = -2636.5
$ws.Range("H46").Value = 975.4828
$ws.Range("I46").Value = 940.4545000000001
$ws.Range("J46").Value = 1085.5714
$ws.Range("K46").Value = 940.4545000000001
$ws.Range("L46").Value = 1085.5714
$ws.Range("M46").Value = -752.4545000000001
$ws.Range("N46").Value = -1461.5714
$ws.Range("H61").Value = 2333.3333
$ws.Range("I61").Value = 2100
$ws.Range("K61").Value = 2100
$ws.Range("M61").Value = -1898
$ws.Range("H113").Value = 2333.3333
$ws.Range("I113").Value = 2100
$ws.Range("K113").Value = 2100
$ws.Range("M113").Value = 70
